$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 originally held only a stray value ("5840917 - Fabricio Maciel Gomes")
# with no label in column A. Delete it; this shifts rows 14-24 up to 13-23,
# matching the new dimension A1:C23.
$ws.Rows.Item(13).Delete()

# Apply the resulting value changes for the shifted rows
$ws.Cells.Item(10, 2).Value = "5840917 - Fabrício Maciel Gomes"
$ws.Cells.Item(10, 3).Value = "5840917 - Fabrício Maciel Gomes"
$ws.Cells.Item(13, 2).Value = "Semestral"
$ws.Cells.Item(13, 3).Value = "Semestral"
$cellR15C2 = $ws.Cells.Item(15, 2)
$cellR15C2.NumberFormat = "@"
$cellR15C2.Value = "01/01/2021"
$cellR15C3 = $ws.Cells.Item(15, 3)
$cellR15C3.NumberFormat = "@"
$cellR15C3.Value = "01/01/2021"
$ws.Cells.Item(18, 2).Value = "5840917 - Fabrício Maciel Gomes"
$ws.Cells.Item(18, 3).Value = "5840917 - Fabrício Maciel Gomes"
$ws.Cells.Item(19, 2).Value = "Aulas expositivas. Trabalhos em grupo. Seminários. Palestras. Exercícios em sala de aula."
$ws.Cells.Item(19, 3).Value = "Aulas expositivas. Trabalhos em grupo. Seminários. Palestras. Exercícios em sala de aula."
$ws.Cells.Item(20, 2).Value = "Duas Provas com peso de 30% cada uma. Trabalhos em sala de aula com peso de 20% e Trabalho final com peso de 20%"
$ws.Cells.Item(20, 3).Value = "Duas Provas com peso de 30% cada uma. Trabalhos em sala de aula com peso de 20% e Trabalho final com peso de 20%"
$ws.Cells.Item(21, 2).Value = "Prova única"
$ws.Cells.Item(21, 3).Value = "Prova única"

# Restore original (General) number formatting/style on the date-like cell(s)
# by pasting formats from an adjacent untouched cell in the same column, since
# directly re-setting NumberFormat would otherwise create a new style entry.
$ws.Cells.Item(14, 2).Copy() | Out-Null
$ws.Cells.Item(15, 2).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Cells.Item(14, 3).Copy() | Out-Null
$ws.Cells.Item(15, 3).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
